$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# "Date: 5/5/2021" -> "Date: 6/05/2021"
#   - the day "5" becomes "6"
#   - a new "0" is inserted right after the following "/" (month becomes "05")
# ---------------------------------------------------------------------
$datePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Date:")) {
        $datePara = $p
        break
    }
}

$labelLen = "Date: ".Length
$pStart = $datePara.Range.Start
$valueRange = $d.Range($pStart + $labelLen, $datePara.Range.End)

# Replace just the first digit (the day) with "6"; wdReplaceOne (1) so only
# the first match in the narrowed range is touched.
$valueRange.Find.Execute("5", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "6", 1) | Out-Null

# Re-acquire the paragraph/value range (length is unchanged so far).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Date:")) {
        $datePara = $p
        break
    }
}
$pStart = $datePara.Range.Start
$valueRange = $d.Range($pStart + $labelLen, $datePara.Range.End)

# Locate the first "/" in the value and insert "0" immediately after it.
$slashRange = $valueRange.Duplicate
$slashRange.Find.Execute("/", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0) | Out-Null
$d.Range($slashRange.End, $slashRange.End).InsertAfter("0")

# ---------------------------------------------------------------------
# "Revision: X02" -> "Revision: X03"
# ---------------------------------------------------------------------
$revPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Revision:")) {
        $revPara = $p
        break
    }
}

$revLabelLen = "Revision: ".Length
$revStart = $revPara.Range.Start
$revValueRange = $d.Range($revStart + $revLabelLen, $revPara.Range.End)
$revValueRange.Find.Execute("2", $true, $false, $false, $false, $false, `
                             $true, 1, $false, "3", 1) | Out-Null
